{"js": "// Applies the cover-letter revision described by the commit:\n// 1) Add a dateline (\"May 22nd, 2023\") above the \"From:\" line.\n// 2) Make the \"From:\" / \"To:\" header use a full blank-line break (like the\n//    rest of the letter) instead of a single line break.\n// 3) Rewrite the body paragraphs (tightened wording, \"job requirements\",\n//    split the long \"co-op program\" paragraph into three paragraphs, and\n//    shorten the closing paragraphs) exactly as in the revised letter.\n\nconst body = context.document.body;\n\n// --- Step 1: insert the dateline before \"From: Austing Dong\" ---\nconst fromSearch = body.search(\"From: Austing Dong\", { matchCase: true });\nfromSearch.load(\"items\");\nawait context.sync();\nif (fromSearch.items.length === 0) {\n  throw new Error('Could not find \"From: Austing Dong\" line');\n}\nfromSearch.items[0].insertText(\"May 22nd, 2023\\v\\v\", Word.InsertLocation.before);\nawait context.sync();\n\n// --- Step 2: turn the single break between \"From:\" and \"To:\" into a full\n//     blank-line break, matching the rest of the letter ---\nconst headerSearch = body.search(\n  \"From: Austing Dong\\vTo: Home Hardware Stores Limited - Divisional Office\",\n  { matchCase: true }\n);\nheaderSearch.load(\"items\");\nawait context.sync();\nif (headerSearch.items.length === 0) {\n  throw new Error(\"Could not find From/To header block\");\n}\nheaderSearch.items[0].insertText(\n  \"From: Austing Dong\\v\\vTo: Home Hardware Stores Limited - Divisional Office\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Step 3: replace the body paragraphs (intro through the closing\n//     \"Thank you ...\" line) with the revised wording/structure ---\nconst oldChunk = \"I am writing to express my strong interest in applying for the position of IT - In Store Systems at Home Hardware Stores Limited - Divisional Office. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirement.\\v\\vI loved computer science as well as developing applications since Middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest, and I did an excellent job in relevant courses in the beginning of my University studies. I found solving business challenges through programming is fascinating because this is the way I feel the sense of accomplishment. Such deep interest in programming and technology has motivated me to deep dive in related fields such as software development, quality assurance and machine learning.\\v\\vThe computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design and data abstraction. Such projects can be viewed on my GitHub: https://github.com/AustingDong. One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information which helps users quickly locate their desired items based on keywords' weight. This application can be used to quickly get all the important items and keywords from NASA Technical Report Server which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project detail can be found here. Through understanding the project requirements, researching on coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment.\\v\\vI am extremely interested in advancing my career and contributing my skills to Home Hardware Stores Limited - Divisional Office. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment.\\v\\vThank you for your valuable time and consideration. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any question or require additional information.\";\nconst newChunk = \"I am writing to express my strong interest in applying for the position of IT - In Store Systems at Home Hardware Stores Limited - Divisional Office. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirements. I would like to highlight the following for your consideration:\\v\\vI have loved computer science and developing applications since middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. My passion for programming continued into my university studies, where I excelled in relevant courses. I find solving business challenges through programming fascinating, as it gives me a sense of accomplishment. This deep interest in programming and technology has motivated me to explore related fields such as software development, quality assurance, and machine learning.\\v\\vThe computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. You can view my projects on my GitHub: https://github.com/AustingDong.\\v\\vOne of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information. This application helps users quickly locate their desired items based on keyword weight and can be used to quickly get all the important items and keywords from NASA Technical Report Server, which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here.\\v\\vThrough understanding project requirements, researching coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment.\\v\\vI am extremely interested in advancing my career and contributing my skills to Home Hardware Stores Limited - Divisional Office. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.\\v\\vThank you for your valuable time and consideration.\";\n\nconst bodySearch = body.search(oldChunk, { matchCase: true });\nbodySearch.load(\"items\");\nawait context.sync();\nif (bodySearch.items.length === 0) {\n  throw new Error(\"Could not find the paragraph block to replace\");\n}\nbodySearch.items[0].insertText(newChunk, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Applies the cover-letter revision described by the commit:\n# 1) Add a dateline (\"May 22nd, 2023\") above the \"From:\" line.\n# 2) Make the \"From:\" / \"To:\" header use a full blank-line break (like the\n#    rest of the letter) instead of a single line break.\n# 3) Rewrite the body paragraphs (tightened wording, \"job requirements\",\n#    split the long \"co-op program\" paragraph into three paragraphs, and\n#    shorten the closing paragraphs) exactly as in the revised letter.\n\n$brk = [char]11\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the dateline before \"From: Austing Dong\" ---\n$range = $d.Content\n$found = $range.Find.Execute(\"From: Austing Dong\")\nif (-not $found) {\n    throw 'Could not find \"From: Austing Dong\" line'\n}\n$range.Collapse(1)\n$range.InsertBefore(\"May 22nd, 2023${brk}${brk}\")\n\n# --- Step 2: turn the single break between \"From:\" and \"To:\" into a full\n#     blank-line break, matching the rest of the letter ---\n$range2 = $d.Content\n$oldHeader = \"From: Austing Dong${brk}To: Home Hardware Stores Limited - Divisional Office\"\n$newHeader = \"From: Austing Dong${brk}${brk}To: Home Hardware Stores Limited - Divisional Office\"\n$found2 = $range2.Find.Execute($oldHeader, $false, $false, $false, $false, $false, $true, 1, $false, $newHeader, 2)\nif (-not $found2) {\n    throw \"Could not find From/To header block\"\n}\n\n# --- Step 3: replace the body paragraphs (intro through the closing\n#     \"Thank you ...\" line) with the revised wording/structure ---\n$range3 = $d.Content\n$oldChunk = \"I am writing to express my strong interest in applying for the position of IT - In Store Systems at Home Hardware Stores Limited - Divisional Office. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirement.${brk}${brk}I loved computer science as well as developing applications since Middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest, and I did an excellent job in relevant courses in the beginning of my University studies. I found solving business challenges through programming is fascinating because this is the way I feel the sense of accomplishment. Such deep interest in programming and technology has motivated me to deep dive in related fields such as software development, quality assurance and machine learning.${brk}${brk}The computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design and data abstraction. Such projects can be viewed on my GitHub: https://github.com/AustingDong. One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information which helps users quickly locate their desired items based on keywords' weight. This application can be used to quickly get all the important items and keywords from NASA Technical Report Server which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project detail can be found here. Through understanding the project requirements, researching on coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment.${brk}${brk}I am extremely interested in advancing my career and contributing my skills to Home Hardware Stores Limited - Divisional Office. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment.${brk}${brk}Thank you for your valuable time and consideration. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any question or require additional information.\"\n$newChunk = \"I am writing to express my strong interest in applying for the position of IT - In Store Systems at Home Hardware Stores Limited - Divisional Office. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirements. I would like to highlight the following for your consideration:${brk}${brk}I have loved computer science and developing applications since middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. My passion for programming continued into my university studies, where I excelled in relevant courses. I find solving business challenges through programming fascinating, as it gives me a sense of accomplishment. This deep interest in programming and technology has motivated me to explore related fields such as software development, quality assurance, and machine learning.${brk}${brk}The computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. You can view my projects on my GitHub: https://github.com/AustingDong.${brk}${brk}One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information. This application helps users quickly locate their desired items based on keyword weight and can be used to quickly get all the important items and keywords from NASA Technical Report Server, which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here.${brk}${brk}Through understanding project requirements, researching coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment.${brk}${brk}I am extremely interested in advancing my career and contributing my skills to Home Hardware Stores Limited - Divisional Office. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.${brk}${brk}Thank you for your valuable time and consideration.\"\n$found3 = $range3.Find.Execute($oldChunk, $false, $false, $false, $false, $false, $true, 1, $false, $newChunk, 2)\nif (-not $found3) {\n    throw \"Could not find the paragraph block to replace\"\n}\n"}
